$d = $word.ActiveDocument

function Replace-One($searchText, $replaceText) {
    $d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null
}

function Replace-All($searchText, $replaceText) {
    $d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

function Replace-After($anchorText, $searchText, $replaceText) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $tail = $d.Range($rng.End, $d.Content.End)
    $tail.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null
}

function Delete-First($searchText) {
    $d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null
}

# Some runs sit immediately after a </w:hyperlink> close. The COM-interop
# engine mis-attributes the hyperlink's run formatting to freshly inserted
# text when the replace/insertion starts right at that boundary (it spans a
# hidden hyperlink boundary marker). Work around it by locating the anchor,
# skipping two positions (the hidden boundary char + the target run's own
# first character), and doing a plain Range.Text assignment for the
# remainder - while keeping that first character literal (unchanged) in
# the replacement text supplied by the caller.
function Replace-AfterHyperlink($anchorText, $keepFirstChar, $restOld, $restNew) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $safeStart = $rng.End + 2
    $tail = $d.Range($safeStart, $d.Content.End)
    $tail.Find.Execute($restOld, $true, $false, $false, $false, $false, $true, 1, $false, $restNew, 1) | Out-Null
}

# 1 & 3: "English" -> "영어" (both occurrences are identical replacements)
Replace-All "English" "영어"

# 2: language list (run starts right after the English hyperlink closes)
Replace-AfterHyperlink "English" " " "/ Portuguese / French / Thai / Vietnamese / Spanish" "/ 포르투갈어 / 프랑스어 / 태국어 / 베트남어 / 스페인어"

# 4: "Brief" -> "간단한 설명"
Replace-One "Brief" "간단한 설명"

# 5: brief description sentence (keep trailing English part intact)
Replace-One "An email sent to partners in the target country who have RSVPed no. It will be sent via customer.io" "RSVP에서 '아니오'로 응답한 초청된 파트너에게 발송되는 이메일입니다. It will be sent via customer.io"

# 6: "Target audience" -> "대상 청중"
Replace-One "Target audience" "대상 청중"

# 7: "Invited partners who RSVP no" -> translated
Replace-One "Invited partners who RSVP no" "RSVP에서 '아니오'로 응답한 초청된 파트너"

# 8: delete the "We'll miss you at the " run entirely
Delete-First "We’ll miss you at the "

# 9: "!" (right after [EVENT NAME], in the headline paragraph) -> translated
Replace-After "[EVENT NAME]" "!" "에서 당신이 그립습니다!"

# 10: delete the "Dear " run entirely
Delete-First "Dear "

# 11: ", " right after [PARTNER NAME] -> " 님, 안녕하세요 "
Replace-After "[PARTNER NAME]" ", " " 님, 안녕하세요 "

# 12: "Thank you for taking the time..." -> "다가오는 "
Replace-One "Thank you for taking the time to respond to our invitation to the upcoming " "다가오는 "

# 13: ". We were really looking forward to seeing you there." -> translated
Replace-One ". We were really looking forward to seeing you there." "에 대한 초대에 응답해 주셔서 감사합니다. 행사장에서 만나 뵙기를 기대했었습니다."

# 14: remove trailing space after "...sometimes come up. "
Replace-One "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up. " "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up."

# 15: "If you’re comfortable sharing..." -> translated
Replace-One "If you’re comfortable sharing it with us, we’d like to know why you responded no. Please reply to this email as your feedback could help us make improvements in our event planning processes and better serve you in the future." "편하시면, 왜 '아니오'로 응답하셨는지 알고 싶습니다. 이 이메일에 회신해 주시면 귀하의 피드백이 향후 행사 기획 프로세스를 개선하고 저희가 더 나은 서비스를 제공하는 데 도움이 될 수 있습니다."

# 16: remove trailing space after "We hope to see you at our future events. "
Replace-One "We hope to see you at our future events. " "We hope to see you at our future events."

# 17: ". " -> "." after WhatsApp hyperlink (end of live-chat/WhatsApp paragraph)
Replace-AfterHyperlink "WhatsApp" "." " " ""

# 18: "If you have any questions, please contact your country manager, " -> "궁금하신 사항은, "
Replace-One "If you have any questions, please contact your country manager, " "궁금하신 사항은, "

# 19: ", at " right after [NAME] -> " 국가 담당자에게 "
Replace-After "[NAME]" ", at " " 국가 담당자에게 "

# 20: " or " right after [EMAIL ADDRESS] -> " 또는 "
Replace-After "[EMAIL ADDRESS]" " or " " 또는 "

# 21: " (WhatsApp). " right after [WHATSAPP NO] -> translated
Replace-After "[WHATSAPP NO]" " (WhatsApp). " " (WhatsApp)으로 연락해 주시기 바랍니다. "

# 22: comment text "choose either one" -> "하나를 선택하세요"
Replace-One "choose either one" "하나를 선택하세요"
